$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)

# Locate the "Content Placeholder 2" shape on this slide (holds the bullet list).
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Content Placeholder 2") {
        $target = $candidate
    }
}

$tr = $target.TextFrame.TextRange

$oldText = "Подобренията на дипломният проект са спазени"
$newPart1 = "Предимства на разработеното в настоящата дипломна работа "
$newPart2 = "уеб приложение"

# Replace the whole old paragraph text with the combined new text first,
# keeping it inside the existing run (preserves rPr / paragraph position).
$fullText = $tr.Text
$startPos = $fullText.IndexOf($oldText) + 1
$oldLen = $oldText.Length
$paraRange = $tr.Characters($startPos, $oldLen)
$paraRange.Text = $newPart1 + $newPart2

# Now re-select just the trailing portion ("уеб приложение") and re-set its
# text in place; this splits the paragraph into two runs, matching the
# authored edit where the text was typed/completed in two passes.
$fullText2 = $tr.Text
$part2Pos = $fullText2.IndexOf($newPart2, $startPos - 1) + 1
$part2Range = $tr.Characters($part2Pos, $newPart2.Length)
$part2Range.Text = $newPart2
